# Updated cryptos list values to match the new OOXML snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.323.42"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "2.289.87"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.31%  "
$ws.Range("D9").Value = "2.285.86"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "2.698.19"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "58.239.66"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "2.323.53"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("B32").Value = "SuiNetwork"
$ws.Range("C32").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.381"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "289.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0952"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0495"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.554"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0211"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("E51").Value = "  +1.35%  "
